# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update row 1 height (76.5 -> 79.5)
$ws.Rows.Item(1).RowHeight = 79.5

# 2) Update title text in B1 (Russian 10.3.1 title), refined wording
$ws.Range("B1").Value = "10.3.1 Доля женщин, сообщивших о том, что в последние 12 месяцев они лично подвергались дискриминации или преследованиям на основаниях, дискриминация по которым запрещена в соответствии с международными стандартами в области прав человека"

# 3) Populate new column E (2023 data) for rows 4-43,
#    copying number format/style from column D then setting value
for ($r = 4; $r -le 43; $r++) {
    $dCell = $ws.Range("D$r")
    $eCell = $ws.Range("E$r")
    $dCell.Copy()
    $eCell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# 4) Set the 2023 values for column E
$ws.Range("E4").Value = 2023
$ws.Range("E5").Value = 6.2
$ws.Range("E7").Value = 7.4
$ws.Range("E8").Value = 5.6
$ws.Range("E10").Value = 4.3
$ws.Range("E11").Value = 7.1
$ws.Range("E12").Value = 2.5
$ws.Range("E13").Value = 2.9
$ws.Range("E14").Value = 3.4
$ws.Range("E15").Value = 1.9
$ws.Range("E16").Value = 9.3000000000000007
$ws.Range("E17").Value = 7.1
$ws.Range("E18").Value = 14.9
$ws.Range("E20").Value = 5.3
$ws.Range("E21").Value = 3.5
$ws.Range("E22").Value = 10
$ws.Range("E23").Value = 5.3
$ws.Range("E24").Value = 5.5
$ws.Range("E25").Value = 7.7
$ws.Range("E26").Value = 6.8
$ws.Range("E27").Value = 5.8
$ws.Range("E28").Value = 7
$ws.Range("E30").Value = "(18,7)"
$ws.Range("E31").Value = 7.5
$ws.Range("E32").Value = 6.1
$ws.Range("E33").Value = 4.9000000000000004
$ws.Range("E34").Value = 6.4
$ws.Range("E36").Value = 32.299999999999997
$ws.Range("E37").Value = 6.1
$ws.Range("E39").Value = 6.7
$ws.Range("E40").Value = 5.7
$ws.Range("E41").Value = 5
$ws.Range("E42").Value = 6.4
$ws.Range("E43").Value = 7.1

# 5) Update footnote row 44 (add ", 2023" to survey years)
$ws.Range("A44").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B44").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C44").Value = "According to the cluster survey for many indicators, 2018, 2023"
